$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.199.92'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -3.74%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.526.61'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -2.95%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '579.04'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '171.95'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -5.36%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.622'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.83%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.516.19'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -3.11%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -6.43%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.69'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +12.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.603'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.66%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '47.42'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -4.59%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.54%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '687.62'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.74%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.092.27'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.97%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '8.83'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.15%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.176.68'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.518.67'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -3.54%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.42%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.46'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -4.57%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.23'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -3.42%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.908'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -3.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '16.61'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -6.66%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '97.87'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -5.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.84'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -4.11%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -6.26%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.44'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -5.50%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '33.34'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -5.41%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.86'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -3.91%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -6.92%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -5.52%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.26'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.31%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '566.42'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -3.26%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.68'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -12.92%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '10.87'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.91%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.93%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '57.24'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -3.76%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.998'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.08%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -2.80%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0441'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -5.63%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.338'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.17%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.449.75'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -6.41%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '33.30'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -6.83%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₃0705'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -7.62%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.49%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -6.90%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.46%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '134.32'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.97%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.17%  '
